# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The "municipio-nombre" column (D) is re-curated from a measure into a
# dimension (refArea / URI-Municipio), matching the pattern already used
# by "provincia-nombre" (E) and "comarca-nombre" (I). The "aragon" column
# (G) is re-curated from a dimension pointing at a local "aragon" mapping
# into a dimension pointing at the shared refArea/Comunidad mapping, and
# its mapping-file cell (row 5) is cleared since it no longer needs its
# own mapping workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("municipio-nombre"): measure -> dimension
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# Column G ("aragon"): local dimension -> shared refArea dimension
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("G4").Value = "URI-Comunidad"

# Row 5 no longer needs a dedicated mapping file for "aragon"
$ws.Range("G5").ClearContents()
